$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New May dates replacing the old April dates, row by row (rows 2..42, step 2)
$dates = @(
    "05/01/20",
    "05/04/20",
    "05/05/20",
    "05/06/20",
    "05/07/20",
    "05/08/20",
    "05/11/20",
    "05/12/20",
    "05/13/20",
    "05/14/20",
    "05/15/20",
    "05/18/20",
    "05/19/20",
    "05/20/20",
    "05/21/20",
    "05/22/20",
    "05/25/20",
    "05/26/20",
    "05/27/20",
    "05/28/20",
    "05/29/20"
)

$row = 2
foreach ($d in $dates) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $d
    $ws.Cells.Item($row, 2).Value = "May"
    $row += 2
}

# Row 44 (previously 04/30/20) is no longer part of the data; clear it
$ws.Cells.Item(44, 1).ClearContents()

$ws.Range("L7").Select()
